$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Technical Challenges",
    "Cybersecurity threats",
    "Lack of User Adoption",
    "Data Security",
    "Change in Requirements",
    "Integration faults",
    "Data Integrity",
    "Scope Creep",
    "Regulatory compliance",
    "Client/Organization Approval",
    "Third Party Approval"
)

$row = 18
foreach ($v in $values) {
    $ws.Range("B$row").Value = $v
    if ($row -le 20) {
        $ws.Rows.Item($row).RowHeight = 15
    }
    $row++
}

$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B29").Select()
